# RN-1289: update the entity associated with a survey resubmission
# Rename the "Facility Code" / "Facility Name" header labels on the
# Test_Facility_Fundamentals sheet to "Entity Code" / "Entity Name".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test_Facility_Fundamentals")

$ws.Range("D2").Value = "Entity Code"
$ws.Range("D3").Value = "Entity Name"
